$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A18").Value = "test"
$ws.Range("C18").Value = "\Testdata\Non_Oncology\DataFiles\LiveSLRPage\NonOnco_Clinical_SLRType_StudyDesign.xlsx"
$ws.Range("B18").Value = "nononco_studydesign_section_validation"
$ws.Range("B18").Select()
